$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-5 from 45233 to 45243
$ws.Range("C2:C5").Value = 45243
